# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Row 16 (Mora 1907): Valor Mora (G16) updated 815000 -> 895000
# - Rows 17-30 (Mora periods for IYOSEIDIS PAEZ JULIO): the "Periodo Mora" (col E)
#   list is reversed into chronological/ascending order (2308 .. 2409 instead of
#   2409 .. 2308), and "Valor Mora" (col G) is updated from 1190000 -> 1543950.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 17 through 30.
$periods = @(
    "2308", "2309", "2310", "2311", "2312",
    "2401", "2402", "2403", "2404", "2405", "2406", "2407", "2408", "2409"
)

$firstRow = 17
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("G$row").Value = 1543950
}

# Row 16 "Valor Mora" update.
$ws.Range("G16").Value = 895000
